# RippleTemplate_Basic.xlsx edit:
#  - add a new "Assay" worksheet (after "Barcodes") holding default
#    assay-transfer settings (Setting / Value columns)
#  - restore the "Patterns" sheet as the active tab with its new selection

$wb = $excel.ActiveWorkbook

# --- add the new "Assay" sheet, right after "Barcodes" -----------------
$barcodes = $wb.Worksheets.Item("Barcodes")
$ws = $wb.Worksheets.Add($null, $barcodes)
$ws.Name = "Assay"

$settings = @(
    @("Setting", "Value"),
    @("DMSO Tolerance", 0.005),
    @("Well Volume (µL)", 25),
    @("Backfill (µL)", 10),
    @("Allowed Error", 0.1),
    @("Destination Replicates", 1),
    @("Use Intermediate Plates", 1),
    @("DMSO Normalization", 1)
)

for ($i = 0; $i -lt $settings.Count; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $settings[$i][0]
    $ws.Cells.Item($row, 2).Value = $settings[$i][1]
}

$ws.Range("A1:B8").Select() | Out-Null

# --- re-activate "Patterns" (the tab that was selected before/after) ---
$patterns = $wb.Worksheets.Item("Patterns")
$patterns.Activate() | Out-Null
$patterns.Range("L17").Select() | Out-Null
